$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the Price column as text (values like "59.956.52" / "0.999" must not
# be auto-converted to numbers by Excel's COM type inference).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '59.956.52'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').Value = '2.414.73'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '552.30'
$ws.Range('E5').Value = '  -0.31%  '
$ws.Range('D6').Value = '137.14'
$ws.Range('E6').Value = '  -1.14%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '0.592'
$ws.Range('E8').Value = '  +3.87%  '
$ws.Range('E9').Value = '  -1.72%  '
$ws.Range('E10').Value = '  -2.83%  '
$ws.Range('E11').Value = '  -0.87%  '
$ws.Range('D12').Value = '0.354'
$ws.Range('E12').Value = '  -1.58%  '
$ws.Range('D13').Value = '25.28'
$ws.Range('E13').Value = '  +1.99%  '
$ws.Range('D14').Value = '2.840.26'
$ws.Range('E14').Value = '  -0.30%  '
$ws.Range('D15').Value = '59.842.86'
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('E16').Value = '  -1.40%  '
$ws.Range('D17').Value = '2.450.63'
$ws.Range('E17').Value = '  +1.80%  '
$ws.Range('D18').Value = '11.33'
$ws.Range('E18').Value = '  -0.80%  '
$ws.Range('D19').Value = '4.42'
$ws.Range('E19').Value = '  -0.08%  '
$ws.Range('D20').Value = '329.01'
$ws.Range('E20').Value = '  -1.33%  '
$ws.Range('D21').Value = '6.68'
$ws.Range('E21').Value = '  -3.20%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').Value = '65.97'
$ws.Range('E23').Value = '  +2.09%  '
$ws.Range('D24').Value = '0.174'
$ws.Range('E24').Value = '  +2.84%  '
$ws.Range('D25').Value = '8.63'
$ws.Range('E25').Value = '  +0.72%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').Value = '1.38'
$ws.Range('E27').Value = '  +0.50%  '
$ws.Range('D28').Value = '0.0₃0772'
$ws.Range('E28').Value = '  -1.79%  '
$ws.Range('E29').Value = '  -2.21%  '
$ws.Range('D30').Value = '169.27'
$ws.Range('E30').Value = '  -0.96%  '
$ws.Range('D31').Value = '6.04'
$ws.Range('E31').Value = '  -3.81%  '
$ws.Range('D32').Value = '18.63'
$ws.Range('E32').Value = '  -0.31%  '
$ws.Range('E33').Value = '  -1.11%  '
$ws.Range('E35').Value = '  -0.55%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('D37').Value = '4.18'
$ws.Range('E37').Value = '  -1.46%  '
$ws.Range('D38').Value = '1.61'
$ws.Range('E38').Value = '  -1.82%  '
$ws.Range('D39').Value = '322.55'
$ws.Range('E39').Value = '  +3.50%  '
$ws.Range('D41').Value = '3.66'
$ws.Range('E41').Value = '  -1.98%  '
$ws.Range('D42').Value = '139.99'
$ws.Range('E42').Value = '  -2.20%  '
$ws.Range('D43').Value = '0.0969'
$ws.Range('E43').Value = '  +0.57%  '
$ws.Range('D44').Value = '19.53'
$ws.Range('E44').Value = '  +2.01%  '
$ws.Range('D45').Value = '0.0515'
$ws.Range('E45').Value = '  -1.77%  '
$ws.Range('D46').Value = '0.579'
$ws.Range('E46').Value = '  +1.32%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = '0.0223'
$ws.Range('E47').Value = '  -1.40%  '
$ws.Range('B48').Value = 'Polygon'
$ws.Range('C48').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D48').Value = '0.387'
$ws.Range('E48').Value = '  -6.11%  '
$ws.Range('D49').Value = '11.04'
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('D50').Value = '1.57'
$ws.Range('E50').Value = '  -3.31%  '
$ws.Range('E51').Value = '  -0.99%  '
